# Quarterly financials update: insert two new quarter columns (D:E) in
# front of the existing data, shifting the historical quarters right,
# and populate the two new columns with the latest reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns at D:E. Everything that was in D:K moves to F:M.
$ws.Columns("D:E").Insert()

# 2) The newly inserted columns come back with the generic/default style.
#    Copy number formatting + styling from column F (which now holds what
#    used to be column D) across the whole table so the new columns look
#    like the columns they were inserted in front of.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the new D/E values (newest two quarters) for every row that
#    carries data. Rows that are entirely blank separators need no values
#    and are intentionally omitted here - the column insert already left
#    them blank.
$data = @{
    7   = @(43465, 43373)
    8   = @(206900, 323800)
    9   = @(5400, 13400)
    10  = @(201500, 310400)
    12  = @("NA", "NA")
    13  = @(0, 0)
    14  = @(0, 0)
    15  = @(118900, 83500)
    17  = @(183500, 138300)
    18  = @(23400, 185500)
    20  = @(-12700, -17400)
    21  = @(289900, 165900)
    22  = @(102200, 82800)
    23  = @(-91500, 85200)
    24  = @(-16400, 28000)
    25  = @(0, 0)
    26  = @(-75100, 57200)
    27  = @(-78900, 53200)
    28  = @(0, 0)
    29  = @(0, "NA")
    30  = @(0, 0)
    31  = @(0, 0)
    32  = @(12700, 17400)
    33  = @(-78900, 53200)
    34  = @(0, 0)
    35  = @(-78900, 53200)
    38  = @(43465, 43373)
    41  = @(631500, 744600)
    42  = @(240800, 145800)
    43  = @(230900, 378900)
    44  = @(18900, 18800)
    45  = @(5500, 9600)
    46  = @(1127700, 1297800)
    47  = @(94500, 92400)
    48  = @(8549200, 8606900)
    49  = @(0, 0)
    50  = @(0, 0)
    51  = @(0, 0)
    52  = @(147600, 175500)
    53  = @(0, 0)
    54  = @(9919000, 10172500)
    57  = @(109400, 82000)
    58  = @(533400, 325300)
    59  = @(96300, 90300)
    60  = @(739100, 497600)
    61  = @(5241800, 5531100)
    62  = @(2182000, 2250500)
    63  = @(0, 0)
    64  = @(0, 0)
    65  = @(0, 0)
    66  = @(8301600, 8414100)
    68  = @(0, 0)
    69  = @(0, 0)
    70  = @(0, 0)
    71  = @(0, 0)
    72  = @(1580700, 1702400)
    73  = @(0, 0)
    74  = @(0, 0)
    75  = @(0, 0)
    76  = @(1617400, 1758500)
    77  = @(0, 0)
    80  = @(43465, 43373)
    81  = @(-78900, 53200)
    83  = @(0, 0)
    84  = @(0, 0)
    85  = @(0, 0)
    86  = @(0, 0)
    87  = @(0, 0)
    88  = @(0, 0)
    89  = @(62700, 175100)
    91  = @("NA", "NA")
    92  = @(0, 0)
    93  = @(0, 0)
    94  = @(-51100, -8300)
    96  = @(-36100, -37000)
    97  = @(0, 0)
    98  = @(0, 0)
    99  = @(0, 0)
    100 = @(-123100, -74500)
    101 = @(-1600, -4900)
    102 = @(-113100, 87400)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}
